$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.774.21'
$ws.Range("E2").Value = '  -3.27%  '
$ws.Range("D3").Value = '1.792.53'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'315.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("D6").Value = "'1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.06%  '
$ws.Range("D7").Value = "'0.5347"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.35%  '
$ws.Range("D8").Value = "'0.3837"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.42%  '
$ws.Range("D9").Value = "'0.07427"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.18%  '
$ws.Range("D10").Value = "'41.44"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.64%  '
$ws.Range("D11").Value = "'1.083"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.00%  '
$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").Value = "'6.206"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.44%  '
$ws.Range("D14").Value = "'7.439"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'20.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.88%  '
$ws.Range("D16").Value = '1.794.25'
$ws.Range("E16").Value = '  -0.54%  '
$ws.Range("D17").Value = "'88.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.52%  '
$ws.Range("D18").Value = "'0.00001058"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.63%  '
$ws.Range("D19").Value = "'0.06526"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("E20").Value = '  +0.02%  '
$ws.Range("D21").Value = "'17.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.24%  '
$ws.Range("D22").Value = "'5.961"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.89%  '
$ws.Range("D23").Value = '27.814.37'
$ws.Range("E23").Value = '  -3.19%  '
$ws.Range("E24").Value = '  -0.51%  '
$ws.Range("D25").Value = "'2.095"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.27%  '
$ws.Range("D26").Value = "'156.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.57%  '
$ws.Range("D27").Value = "'20.21"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.21%  '
$ws.Range("D28").Value = '1.999.96'
$ws.Range("E28").Value = '  -0.71%  '
$ws.Range("D29").Value = "'2.327"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.55%  '
$ws.Range("D30").Value = "'121.32"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.49%  '
$ws.Range("D31").Value = "'1.111"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.71%  '
$ws.Range("E32").Value = '  +3.60%  '
$ws.Range("D33").Value = "'3.653"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.02%  '
$ws.Range("D34").Value = "'5.503"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.68%  '
$ws.Range("D35").Value = "'0.06977"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +8.49%  '
$ws.Range("D36").Value = "'0.2191"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.91%  '
$ws.Range("D37").Value = "'0.02273"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.63%  '
$ws.Range("D38").Value = "'5.052"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").Value = "'11.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.70%  '
$ws.Range("D40").Value = "'8.402"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.45%  '
$ws.Range("D41").Value = "'0.6098"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'1.411"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.02%  '
$ws.Range("D43").Value = "'1.159"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.79%  '
$ws.Range("D44").Value = "'13.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("D45").Value = "'3.678"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.45%  '
$ws.Range("D46").Value = "'0.5699"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.88%  '
$ws.Range("D47").Value = "'125.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.70%  '
$ws.Range("D48").Value = "'1.907"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.26%  '
$ws.Range("D49").Value = "'1.169"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.57%  '
$ws.Range("D50").Value = "'0.06786"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.52%  '
$ws.Range("D51").Value = "'71.73"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.31%  '
